$d = $word.ActiveDocument

$rng = $d.Content
$found = $rng.Find.Execute("set_category(category)")
$para = $rng.Paragraphs(1)
$origStart = $para.Range.Start

$collapsedRange = $d.Range($origStart, $origStart)
Write-Output "collapsed range text: [$($collapsedRange.Text)]"

$xmlFragment = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$collapsedRange.InsertXML($xmlFragment)
